$d = $word.ActiveDocument

# --- Programa (Portuguese) paragraph: insert manual line breaks ---
$d.Content.Find.Execute(
    "desigualdades. Funções Reais",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "desigualdades. ^lFunções Reais", 2) | Out-Null

$d.Content.Find.Execute(
    "funções hiperbólicas. Modelagem",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "funções hiperbólicas. ^lModelagem", 2) | Out-Null

# --- Programa (English / italic) paragraph: insert manual line breaks ---
$d.Content.Find.Execute(
    "inequalities.Real Functions",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "inequalities.^lReal Functions", 2) | Out-Null

$d.Content.Find.Execute(
    "hyperbolic functions.Modeling",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "hyperbolic functions.^lModeling", 2) | Out-Null

# --- Bibliografia paragraph: insert double manual line breaks between entries ---
$d.Content.Find.Execute(
    "v.1.ANTON, Howard",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "v.1.^l^lANTON, Howard", 2) | Out-Null

$d.Content.Find.Execute(
    "2007.THOMAS, George",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "2007.^l^lTHOMAS, George", 2) | Out-Null

$d.Content.Find.Execute(
    "v.1,FLEMMING, Diva",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "v.1,^l^lFLEMMING, Diva", 2) | Out-Null
